# Auto-generated edit script applying the Pandaemonium_Profits diff
# Updates plain numeric cell values (and a few cell additions/removals) across
# the ALC, ARM, BSM, CRP, CUL, GSM, and LTW sheets.

$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H19").Value = 55556676
$wsALC.Range("I19").Value = 166666980
$wsALC.Range("K19").Value = 166666980
$wsALC.Range("M19").Value = -166666805
$wsALC.Range("H32").Value = 422.75
$wsALC.Range("I32").Value = 321.125
$wsALC.Range("K32").Value = 321.125
$wsALC.Range("M32").Value = 4.875
$wsALC.Range("H33").Value = 101.69231
$wsALC.Range("I33").Value = 64
$wsALC.Range("K33").Value = 64
$wsALC.Range("M33").Value = 165
$wsALC.Range("H98").Value = 1590.4445
$wsALC.Range("I98").Value = 1308.3158
$wsALC.Range("J98").Value = 2260.5
$wsALC.Range("K98").Value = 1308.3158
$wsALC.Range("L98").Value = 2260.5
$wsALC.Range("M98").Value = 189.6841999999999
$wsALC.Range("N98").Value = -5256.5
$wsALC.Range("H100").Value = 1502.16
$wsALC.Range("I100").Value = 1313.1111
$wsALC.Range("J100").Value = 1988.2858
$wsALC.Range("K100").Value = 1313.1111
$wsALC.Range("L100").Value = 1988.2858
$wsALC.Range("M100").Value = -772.1111000000001
$wsALC.Range("N100").Value = -3070.2858
$wsALC.Range("H111").Value = 3329.3333
$wsALC.Range("I111").Value = 3433.3333
$wsALC.Range("J111").Value = 3277.3333
$wsALC.Range("K111").Value = 10299.9999
$wsALC.Range("L111").Value = 9831.999899999999
$wsALC.Range("M111").Value = -7232.999899999999
$wsALC.Range("N111").Value = -15965.9999
$wsALC.Range("H122").Value = 1590.4445
$wsALC.Range("I122").Value = 1308.3158
$wsALC.Range("J122").Value = 2260.5
$wsALC.Range("K122").Value = 3924.9474
$wsALC.Range("L122").Value = 6781.5
$wsALC.Range("M122").Value = -1474.9474
$wsALC.Range("N122").Value = -11681.5
$wsALC.Range("H134").Value = 40495
$wsALC.Range("J134").Value = 40495
$wsALC.Range("L134").Value = 40495
$wsALC.Range("N134").Value = -50635

$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("H32").Value = 17667.953
$wsARM.Range("I32").Value = 19660.893
$wsARM.Range("J32").Value = 5267.4443
$wsARM.Range("K32").Value = 19660.893
$wsARM.Range("L32").Value = 5267.4443
$wsARM.Range("M32").Value = -19373.893
$wsARM.Range("N32").Value = -5841.4443
$wsARM.Range("H61").Value = 7787.9287
$wsARM.Range("I61").Value = 4404.8823
$wsARM.Range("J61").Value = 22165.875
$wsARM.Range("K61").Value = 4404.8823
$wsARM.Range("L61").Value = 22165.875
$wsARM.Range("M61").Value = -4192.8823
$wsARM.Range("N61").Value = -22589.875
$wsARM.Range("H62").Value = 0
$wsARM.Range("J62").Value = 0
$wsARM.Range("L62").Value = 0
$wsARM.Range("N62").ClearContents()
$wsARM.Range("H63").Value = 3428.6
$wsARM.Range("I63").Value = 2755.2856
$wsARM.Range("J63").Value = 4999.6665
$wsARM.Range("K63").Value = 2755.2856
$wsARM.Range("L63").Value = 4999.6665
$wsARM.Range("M63").Value = -2069.2856
$wsARM.Range("N63").Value = -6371.6665
$wsARM.Range("H64").Value = 40091
$wsARM.Range("J64").Value = 40091
$wsARM.Range("L64").Value = 40091
$wsARM.Range("N64").Value = -40587
$wsARM.Range("H65").Value = 0
$wsARM.Range("J65").Value = 0
$wsARM.Range("L65").Value = 0
$wsARM.Range("N65").ClearContents()
$wsARM.Range("H66").Value = 3428.6
$wsARM.Range("I66").Value = 2755.2856
$wsARM.Range("J66").Value = 4999.6665
$wsARM.Range("K66").Value = 13776.428
$wsARM.Range("L66").Value = 24998.3325
$wsARM.Range("M66").Value = -10344.428
$wsARM.Range("N66").Value = -31862.3325
$wsARM.Range("H67").Value = 40091
$wsARM.Range("J67").Value = 40091
$wsARM.Range("L67").Value = 40091
$wsARM.Range("N67").Value = -41807
$wsARM.Range("H68").Value = 0
$wsARM.Range("J68").Value = 0
$wsARM.Range("L68").Value = 0
$wsARM.Range("N68").ClearContents()
$wsARM.Range("H71").Value = 0
$wsARM.Range("J71").Value = 0
$wsARM.Range("L71").Value = 0
$wsARM.Range("N71").ClearContents()
$wsARM.Range("H76").Value = 0
$wsARM.Range("I76").Value = 0
$wsARM.Range("K76").Value = 0
$wsARM.Range("M76").ClearContents()
$wsARM.Range("H79").Value = 0
$wsARM.Range("I79").Value = 0
$wsARM.Range("K79").Value = 0
$wsARM.Range("M79").ClearContents()
$wsARM.Range("H125").Value = 0
$wsARM.Range("J125").Value = 0
$wsARM.Range("L125").Value = 0
$wsARM.Range("N125").ClearContents()
$wsARM.Range("H132").Value = 1735.1428
$wsARM.Range("I132").Value = 1502.4286
$wsARM.Range("J132").Value = 2200.5715
$wsARM.Range("K132").Value = 4507.2858
$wsARM.Range("L132").Value = 6601.7145
$wsARM.Range("M132").Value = -1977.2858
$wsARM.Range("N132").Value = -11661.7145
$wsARM.Range("H136").Value = 7787.9287
$wsARM.Range("I136").Value = 4404.8823
$wsARM.Range("J136").Value = 22165.875
$wsARM.Range("K136").Value = 13214.6469
$wsARM.Range("L136").Value = 66497.625
$wsARM.Range("M136").Value = -10664.6469
$wsARM.Range("N136").Value = -71597.625

$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("H22").Value = 234.75
$wsBSM.Range("I22").Value = 234.75
$wsBSM.Range("K22").Value = 234.75
$wsBSM.Range("M22").Value = -61.75
$wsBSM.Range("H94").Value = 1243.9445
$wsBSM.Range("I94").Value = 1009.6
$wsBSM.Range("J94").Value = 1536.875
$wsBSM.Range("K94").Value = 1009.6
$wsBSM.Range("L94").Value = 1536.875
$wsBSM.Range("M94").Value = -558.6
$wsBSM.Range("N94").Value = -2438.875
$wsBSM.Range("H107").Value = 1762.2
$wsBSM.Range("I107").Value = 905.5
$wsBSM.Range("K107").Value = 905.5
$wsBSM.Range("M107").Value = 1014.5
$wsBSM.Range("H134").Value = 30723.514
$wsBSM.Range("I134").Value = 2245.84
$wsBSM.Range("J134").Value = 101917.7
$wsBSM.Range("K134").Value = 6737.52
$wsBSM.Range("L134").Value = 305753.1
$wsBSM.Range("M134").Value = -4202.52
$wsBSM.Range("N134").Value = -310823.1

$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Range("H31").Value = 1683.2222
$wsCRP.Range("I31").Value = 1415.0741
$wsCRP.Range("J31").Value = 2487.6667
$wsCRP.Range("K31").Value = 1415.0741
$wsCRP.Range("L31").Value = 2487.6667
$wsCRP.Range("M31").Value = -1120.0741
$wsCRP.Range("N31").Value = -3077.6667
$wsCRP.Range("H34").Value = 1683.2222
$wsCRP.Range("I34").Value = 1415.0741
$wsCRP.Range("J34").Value = 2487.6667
$wsCRP.Range("K34").Value = 1415.0741
$wsCRP.Range("L34").Value = 2487.6667
$wsCRP.Range("M34").Value = -1213.0741
$wsCRP.Range("N34").Value = -2891.6667
$wsCRP.Range("H58").Value = 2458927.2
$wsCRP.Range("I58").Value = 3789350.5
$wsCRP.Range("J58").Value = 2761.6924
$wsCRP.Range("K58").Value = 3789350.5
$wsCRP.Range("L58").Value = 2761.6924
$wsCRP.Range("M58").Value = -3789147.5
$wsCRP.Range("N58").Value = -3167.6924
$wsCRP.Range("H94").Value = 961.64703
$wsCRP.Range("I94").Value = 1212
$wsCRP.Range("J94").Value = 946
$wsCRP.Range("K94").Value = 1212
$wsCRP.Range("L94").Value = 946
$wsCRP.Range("M94").Value = -761
$wsCRP.Range("N94").Value = -1848
$wsCRP.Range("H134").Value = 2098.95
$wsCRP.Range("I134").Value = 1885.6364
$wsCRP.Range("J134").Value = 3104.5715
$wsCRP.Range("K134").Value = 5656.9092
$wsCRP.Range("L134").Value = 9313.7145
$wsCRP.Range("M134").Value = -3121.9092
$wsCRP.Range("N134").Value = -14383.7145
$wsCRP.Range("H136").Value = 2458927.2
$wsCRP.Range("I136").Value = 3789350.5
$wsCRP.Range("J136").Value = 2761.6924
$wsCRP.Range("K136").Value = 11368051.5
$wsCRP.Range("L136").Value = 8285.0772
$wsCRP.Range("M136").Value = -11365501.5
$wsCRP.Range("N136").Value = -13385.0772

$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("H5").Value = 18530126
$wsCUL.Range("I5").Value = 567.8
$wsCUL.Range("K5").Value = 1703.4
$wsCUL.Range("M5").Value = -1591.4
$wsCUL.Range("H68").Value = 7354.6113
$wsCUL.Range("I68").Value = 800
$wsCUL.Range("J68").Value = 8665.532999999999
$wsCUL.Range("K68").Value = 2400
$wsCUL.Range("L68").Value = 25996.599
$wsCUL.Range("M68").Value = -1589
$wsCUL.Range("N68").Value = -27618.599
$wsCUL.Range("H71").Value = 7354.6113
$wsCUL.Range("I71").Value = 800
$wsCUL.Range("J71").Value = 8665.532999999999
$wsCUL.Range("K71").Value = 7200
$wsCUL.Range("L71").Value = 77989.79699999999
$wsCUL.Range("M71").Value = -3144
$wsCUL.Range("N71").Value = -86101.79699999999
$wsCUL.Range("H124").Value = 2905
$wsCUL.Range("J124").Value = 3630.8333
$wsCUL.Range("L124").Value = 10892.4999
$wsCUL.Range("N124").Value = -20712.4999
$wsCUL.Range("H135").Value = 18530126
$wsCUL.Range("I135").Value = 567.8
$wsCUL.Range("K135").Value = 5110.2
$wsCUL.Range("M135").Value = -2575.2

$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("H102").Value = 4007.9
$wsGSM.Range("I102").Value = 3562.3809
$wsGSM.Range("K102").Value = 3562.3809
$wsGSM.Range("M102").Value = -1940.3809
$wsGSM.Range("H103").Value = 10000
$wsGSM.Range("J103").Value = 10000
$wsGSM.Range("L103").Value = 10000
$wsGSM.Range("N103").Value = -12344
$wsGSM.Range("H122").Value = 7544.3335
$wsGSM.Range("I122").Value = 11399.8
$wsGSM.Range("K122").Value = 34199.39999999999
$wsGSM.Range("M122").Value = -31749.39999999999

$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("H40").Value = 4069.625
$wsLTW.Range("I40").Value = 3761.6
$wsLTW.Range("K40").Value = 3761.6
$wsLTW.Range("M40").Value = -3625.6
$wsLTW.Range("H55").Value = 435.15384
$wsLTW.Range("I55").Value = 380.83334
$wsLTW.Range("J55").Value = 481.7143
$wsLTW.Range("K55").Value = 380.83334
$wsLTW.Range("L55").Value = 481.7143
$wsLTW.Range("M55").Value = -207.83334
$wsLTW.Range("N55").Value = -827.7143
$wsLTW.Range("H61").Value = 1277791.6
$wsLTW.Range("I61").Value = 42467.6
$wsLTW.Range("J61").Value = 3336665
$wsLTW.Range("K61").Value = 42467.6
$wsLTW.Range("L61").Value = 3336665
$wsLTW.Range("M61").Value = -42265.6
$wsLTW.Range("N61").Value = -3337069
$wsLTW.Range("H113").Value = 1277791.6
$wsLTW.Range("I113").Value = 42467.6
$wsLTW.Range("J113").Value = 3336665
$wsLTW.Range("K113").Value = 42467.6
$wsLTW.Range("L113").Value = 3336665
$wsLTW.Range("M113").Value = -40297.6
$wsLTW.Range("N113").Value = -3341005
$wsLTW.Range("H122").Value = 6076.615
$wsLTW.Range("I122").Value = 5241.5527
$wsLTW.Range("J122").Value = 8343.214
$wsLTW.Range("K122").Value = 15724.6581
$wsLTW.Range("L122").Value = 25029.642
$wsLTW.Range("M122").Value = -13274.6581
$wsLTW.Range("N122").Value = -29929.642
